$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tipos de Datos")

$ws.Range("A7").Value = "DESC"
$ws.Range("B7").Value = "STRING"
$ws.Range("C7").Value = 32
$ws.Range("D7").Value = "Descripcion De respuestas"

$ws.Range("A8").Select()
